$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price (D) column cells we touch keep storing their numeric-looking
# values as text (matching the source data format) instead of being
# auto-converted to Excel numbers, which would lose exact trailing-zero
# formatting (e.g. "0.09320" vs 0.0932).
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D23", "D24", "D25", "D26", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "272.38"
$ws.Range("D3").Value = "23.08"
$ws.Range("D4").Value = "6.372"
$ws.Range("D5").Value = "0.06280"
$ws.Range("D6").Value = "3.650"
$ws.Range("D7").Value = "6.724"
$ws.Range("D8").Value = "1.388"
$ws.Range("D9").Value = "0.8393"
$ws.Range("D10").Value = "0.1630"
$ws.Range("D11").Value = "0.08465"
$ws.Range("D12").Value = "0.03471"
$ws.Range("D13").Value = "0.03137"
$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D14").Value = "3.997"
$ws.Range("E14").Value = "13MCDexMCB"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "0.09320"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("D16").Value = "0.001730"
$ws.Range("D17").Value = "0.04854"
$ws.Range("D18").Value = "0.006236"
$ws.Range("D19").Value = "0.005490"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("D20").Value = "0.001089"
$ws.Range("D21").Value = "0.0001500"
$ws.Range("D23").Value = "2.315"
$ws.Range("D24").Value = "0.01383"
$ws.Range("D25").Value = "0.3406"
$ws.Range("D26").Value = "0.1262"
$ws.Range("D40").Value = "0.04692"
$ws.Range("D41").Value = "0.006894"
$ws.Range("D42").Value = "0.1177"
$ws.Range("D43").Value = "0.003456"
$ws.Range("D44").Value = "0.01260"
$ws.Range("D45").Value = "0.00006247"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("D47").Value = "0.7971"
$ws.Range("D48").Value = "0.09689"
$ws.Range("D49").Value = "0.00001400"
$ws.Range("E49").Value = "48CryptobidCoinCBCWorstin24h"
$ws.Range("D50").Value = "0.01240"
